$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update the Files query text in B4: remove the "File Type" and "Breed" lines
# and adjust indentation of the following lines, matching the corrected
# ICDC Breed 1-14 scripts.
$newFilesQuery = "MATCH (f:file)-->(parent)`nWITH DISTINCT f, parent`nMATCH (f)-[*]->(c:case)<--(demo:demographic)`nWHERE demo.breed IN ['Irish Setter']`nOPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)`nOPTIONAL MATCH (samp:sample)-->(c)`nWITH DISTINCT f, parent, c, demo, diag, s`nRETURN  coalesce(f.file_name, '') AS ``File Name``,`n           coalesce(labels(parent)[0], '') AS ``Association``,`n        coalesce(f.file_description, '') AS ``Description``,`n        coalesce(f.file_format, '') AS ``Format``,`n        coalesce(f.file_size, '') AS ``Size``,`n        coalesce(c.case_id, '') AS ``Case ID``,`n         coalesce(diag.disease_term,'') AS Diagnosis , `n        coalesce(s.clinical_study_designation,'') AS ``Study Code``"

$ws.Range("B4").Value = $newFilesQuery

# The removed lines shorten the wrapped text, so the row shrinks accordingly
# (was 17 wrapped lines / 246.5pt, now 15 wrapped lines / 217.5pt).
$ws.Rows.Item(4).RowHeight = 217.5

# Reflect the navigation/selection state after the edit (user ended up viewing B4)
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("B4").Select()
